# Updates the cryptos list: refreshed Price (column D) and Volume(1h) (column E)
# figures, plus a swap of the ImmutableX / EthereumClassic rows (33 <-> 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a value into a cell while forcing it to stay TEXT ------
# Many "Price" values look like plain numbers (e.g. "1.00", "7.22") and a
# bare .Value assignment would get auto-coerced into a numeric cell, which
# does not match how this sheet stores its data (plain text, same as the
# non-numeric-looking prices like "67.440.75"). Briefly marking the cell as
# Text (@) before the assignment keeps it a string, and restoring the style
# afterwards leaves the cell's own formatting untouched.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Column D ("Price") updates
# ---------------------------------------------------------------------
$dValues = [ordered]@{
    "D2"  = "67.440.75"
    "D3"  = "3.757.88"
    "D4"  = "1.00"
    "D5"  = "595.22"
    "D6"  = "169.94"
    "D7"  = "3.755.51"
    "D14" = "36.69"
    "D15" = "4.388.37"
    "D17" = "18.71"
    "D18" = "67.535.75"
    "D19" = "7.22"
    "D21" = "10.53"
    "D22" = "469.84"
    "D23" = "0.722"
    "D25" = "83.85"
    "D26" = "2.23"
    "D27" = "12.18"
    "D28" = "10.37"
    "D29" = "0.999"
    "D31" = "3.909.85"
    "D32" = "7.70"
    "D35" = "9.15"
    "D36" = "3.721.47"
    "D39" = "0.139"
    "D41" = "0.997"
    "D42" = "1.00"
    "D43" = "0.312"
    "D45" = "8.75"
    "D47" = "45.91"
    "D48" = "398.29"
    "D50" = "141.33"
}
foreach ($addr in $dValues.Keys) {
    Set-TextValue $addr $dValues[$addr]
}

# ---------------------------------------------------------------------
# Row 33 / 34 swap: ImmutableX and EthereumClassic trade places (name,
# link and price move; each row's own Volume(1h) cell stays put).
# ---------------------------------------------------------------------
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "30.51"

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D34" "2.24"

# ---------------------------------------------------------------------
# Column E ("Volume(1h)") updates -- every row except 33, whose value is
# unchanged by this edit.
# ---------------------------------------------------------------------
$eValues = [ordered]@{
    "E2"  = "  -0.71%  "
    "E3"  = "  -2.06%  "
    "E4"  = "  +0.04%  "
    "E5"  = "  -0.84%  "
    "E6"  = "  +1.22%  "
    "E7"  = "  -2.14%  "
    "E8"  = "  +0.00%  "
    "E9"  = "  -0.20%  "
    "E10" = "  +1.14%  "
    "E11" = "  +0.04%  "
    "E12" = "  -0.57%  "
    "E13" = "  +7.33%  "
    "E14" = "  -1.11%  "
    "E15" = "  -2.03%  "
    "E16" = "  -1.86%  "
    "E17" = "  +2.08%  "
    "E18" = "  -0.59%  "
    "E19" = "  -2.48%  "
    "E20" = "  +0.97%  "
    "E21" = "  -5.19%  "
    "E22" = "  +0.80%  "
    "E23" = "  -1.76%  "
    "E24" = "  -8.73%  "
    "E25" = "  +1.14%  "
    "E26" = "  -0.03%  "
    "E27" = "  +0.52%  "
    "E28" = "  +3.24%  "
    "E29" = "  -0.15%  "
    "E30" = "  -1.51%  "
    "E31" = "  -1.88%  "
    "E32" = "  +0.72%  "
    "E34" = "  -2.95%  "
    "E35" = "  -4.30%  "
    "E36" = "  -2.03%  "
    "E37" = "  +5.02%  "
    "E38" = "  +0.86%  "
    "E39" = "  -0.99%  "
    "E40" = "  -0.62%  "
    "E41" = "  -1.95%  "
    "E42" = "  +0.00%  "
    "E43" = "  -0.65%  "
    "E44" = "  +0.01%  "
    "E45" = "  +0.12%  "
    "E46" = "  -1.13%  "
    "E47" = "  -2.09%  "
    "E48" = "  -5.18%  "
    "E49" = "  -8.14%  "
    "E50" = "  -0.79%  "
    "E51" = "  -0.38%  "
}
foreach ($addr in $eValues.Keys) {
    $ws.Range($addr).Value = $eValues[$addr]
}
